# Fix code that prevents figures overlapping
# Clear the figure-number sequence values in C10:I10 (keep formatting),
# and update the active selection to I10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("control_panel")

# Clear the contents (values) of C10:I10 while preserving cell formatting.
$ws.Range("C10:I10").ClearContents()

# Update the selected / active cell on the sheet to I10.
$ws.Activate()
$ws.Range("I10").Select()
